$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 126-127, pushing the former rows 126-127 down to 128-129.
$ws.Rows("126:127").Insert()

# Copy the index-column formatting (style) down onto the two freshly inserted rows
# so A126/A127 keep the same bold/bordered/centered style as the rest of column A.
$ws.Range("A125").Copy()
$ws.Range("A126:A127").PasteSpecial(-4122)

# New row 126 (index 125): Roman Ramirez (G) and Will Simpson (H) played.
$ws.Range("A126").Value = 125
$ws.Range("B126:F126").Value = "'"
$ws.Range("G126").Value = 1200.544650217021
$ws.Range("H126").Value = 1207.406863229328

# New row 127 (index 126): Roman Ramirez (G) and Will Simpson (H) played again.
$ws.Range("A127").Value = 126
$ws.Range("B127:F127").Value = "'"
$ws.Range("G127").Value = 1185.628639295045
$ws.Range("H127").Value = 1222.322874151305

# The "empty" cells above were typed with a leading apostrophe just to force an
# empty *text* cell (matching the rest of the sheet's blank-cell convention)
# instead of Excel dropping them entirely; strip the resulting quote-prefix
# formatting back to the plain, unstyled look used everywhere else.
$ws.Range("B2").Copy()
$ws.Range("B126:F127").PasteSpecial(-4122)

# The two rows that were pushed down (former rows 126-127) keep their own data but
# their running index in column A must advance by 2 to stay sequential.
$ws.Range("A128").Value = 127
$ws.Range("A129").Value = 128

# The old row 127 (now shifted to row 129) gets an updated Will Simpson (H) Elo value.
$ws.Range("H129").Value = 1234.740270956793
